# Feature: add arrows (arrow_n). Fixed bugs, removed unnecessary code.
#
# The "meta" sheet stores chart configuration as key/value rows in columns
# A/B, terminated by a single blank (but still styled) divider row 13.
# Insert a new row above that divider so it becomes row 14, then fill the
# freed-up row 13 with the new "style" = "default" setting. The inserted
# row inherits the key column's bold/orange formatting from row 12, same
# as the divider row had.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "style"
$ws.Range("B13").Value = "default"
